$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").ClearContents()
$ws.Range("C2").Value = 0
